# Correção ao Diagrama Lógico e Use case de Ementa Semanal
#
# The original sheet had a spare/empty row (row 9) inside the "Cenário
# Normal" block (merged B6:B10) that didn't belong there. Removing it
# shifts every row below up by one (rows 10-18 become 9-17), shrinks the
# used range, and shifts the dependent merged cells accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete the stray row 9 - Excel shifts rows 10:18 up to 9:17 and adjusts
# the B6:B10 / B11:B14 / B15:B18 merges to B6:B9 / B10:B13 / B14:B17.
$ws.Rows(9).Delete() | Out-Null

# The author ended up with the selection on D22 (previously C22).
$ws.Range("D22").Select() | Out-Null
